$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# ---- Row 23 ----
$ws.Range("A23").Value = 111664006
$ws.Range("B23").Value = 96368
Set-TextValue $ws.Range("D23") 'LC'
$ws.Range("E23").Value = 221952
Set-TextValue $ws.Range("F23") 'Spindelblomster'
Set-TextValue $ws.Range("G23") 'Neottia cordata'
Set-TextValue $ws.Range("H23") '(L.) Rich.'
Set-TextValue $ws.Range("I23") '25'
Set-TextValue $ws.Range("J23") 'stjälkar/strån/skott'
Set-TextValue $ws.Range("K23") 'överblommad'
$ws.Range("L23").ClearContents() | Out-Null
Set-TextValue $ws.Range("N23") 'observerad'
$ws.Range("Q23").Value = 485633.0508789647
$ws.Range("R23").Value = 7005896.397059018
$ws.Range("S23").Value = 10
Set-TextValue $ws.Range("AC23") 'Minst 25 st. skott/stjälkar av spindelblomster varav minst ett skott som bär en överblommad blomstängel. Obs! Fyndplatsen ligger inom en avverkningsanmäld yta med beteckning A 32699-2023.'
$ws.Range("AJ23").ClearContents() | Out-Null
$ws.Range("AK23").ClearContents() | Out-Null
$ws.Range("AO23").ClearContents() | Out-Null
# ---- Row 24 ----
$ws.Range("A24").Value = 111663296
$ws.Range("B24").Value = 77515
Set-TextValue $ws.Range("D24") 'NT'
$ws.Range("E24").Value = 6425
Set-TextValue $ws.Range("F24") 'Garnlav'
Set-TextValue $ws.Range("G24") 'Alectoria sarmentosa'
Set-TextValue $ws.Range("H24") '(Ach.) Ach.'
$ws.Range("I24").ClearContents() | Out-Null
$ws.Range("J24").ClearContents() | Out-Null
$ws.Range("K24").ClearContents() | Out-Null
$ws.Range("L24").ClearContents() | Out-Null
$ws.Range("N24").ClearContents() | Out-Null
$ws.Range("Q24").Value = 485609.6900141542
$ws.Range("R24").Value = 7005829.216201009
Set-TextValue $ws.Range("AC24") 'Enstaka bålar av garnlav på gran. Obs! Fyndplatsen ligger inom en avverkningsanmäld yta med beteckning A 32699-2023.'
Set-TextValue $ws.Range("AJ24") 'gran'
Set-TextValue $ws.Range("AK24") 'Picea abies'
Set-TextValue $ws.Range("AO24") 'Picea abies'
# ---- Row 25 ----
$ws.Range("A25").Value = 111662960
$ws.Range("B25").Value = 96348
Set-TextValue $ws.Range("D25") 'VU'
$ws.Range("E25").Value = 220787
Set-TextValue $ws.Range("F25") 'Knärot'
Set-TextValue $ws.Range("G25") 'Goodyera repens'
Set-TextValue $ws.Range("H25") '(L.) R. Br.'
Set-TextValue $ws.Range("I25") '880'
$ws.Range("Q25").Value = 485592.9227098347
$ws.Range("R25").Value = 7005821.162446524
$ws.Range("S25").Value = 5
Set-TextValue $ws.Range("AC25") 'Här finns rikligt med knärot i en mer luckig del av skogen. Minst 880 st. skott/stjälkar och 63 st. överblommade blomstänglar av knärot inom en yta på ca 6 m2. Obs! Fyndplatsen ligger inom en avverkningsanmäld yta med beteckning A 32699-2023.'
